$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 50, shifting existing rows 50:107 down to 51:108
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly record
$ws.Cells.Item(50, 1).Value = 10
$ws.Cells.Item(50, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(50, 3).Value = "La Araucanía"
$ws.Cells.Item(50, 4).Value = 45128
$ws.Cells.Item(50, 5).Value = 9
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100108
$ws.Cells.Item(50, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(50, 9).Value = 100108007
$ws.Cells.Item(50, 10).Value = "Coco"
$ws.Cells.Item(50, 11).Value = "Sin especificar"
$ws.Cells.Item(50, 12).Value = "Primera"
$ws.Cells.Item(50, 13).Value = 25
$ws.Cells.Item(50, 14).Value = 36000
$ws.Cells.Item(50, 15).Value = 36000
$ws.Cells.Item(50, 16).Value = 36000
$ws.Cells.Item(50, 17).Value = "`$/malla 20 unidades"
$ws.Cells.Item(50, 18).Value = "Perú"
$ws.Cells.Item(50, 19).Value = 1800
$ws.Cells.Item(50, 20).Value = 20
